{"js": "// Replace the Changeset number \"49473\" with \"49511\" in the table cell\n// (Runbook \"Get the following from TFS...\" step). The original text is\n// split across several runs (\"49\" / \"4\" / \"73\"). Only replace the \"473\"\n// portion (which spans the last two runs) so the leading \"49\" run is\n// left untouched, matching the authored edit.\nconst body = context.document.body;\nconst results = body.search(\"473\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find target text '473' to replace.\");\n}\n\nresults.items[0].insertText(\"511\", \"Replace\");\nawait context.sync();\n", "ps1": "# Update the Changeset number from 49473 to 49511 in the Runbook step\n# that reads \"... Changeset 49473\" (Step 1 of the \"Get the zip from TFS\"\n# procedure). Only the \"473\" substring is targeted so the leading \"49\"\n# text is left in place and the edit stays scoped to the changeset\n# number itself.\n$d = $word.ActiveDocument\n\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.Text = \"473\"\n$rng.Find.Replacement.ClearFormatting()\n$rng.Find.Replacement.Text = \"511\"\n$rng.Find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, \"511\", 2) | Out-Null\n"}
